$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 699, shifting rows 699:729 down
# to 701:731 (matches the target: dimension grows from A1:R729 to A1:R731).
$ws.Rows.Item(699).Insert()
$ws.Rows.Item(699).Insert()

# Populate new row 699
$ws.Cells.Item(699, 1).Value = 9
$ws.Cells.Item(699, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(699, 3).Value = "Metropolitana"
$ws.Cells.Item(699, 4).Value = 44509
$ws.Cells.Item(699, 5).Value = 13
$ws.Cells.Item(699, 6).Value = 100112045
$ws.Cells.Item(699, 7).Value = "Zapallo"
$ws.Cells.Item(699, 8).Value = "Camote"
$ws.Cells.Item(699, 9).Value = "1a nueva(o)"
$ws.Cells.Item(699, 10).Value = 133
$ws.Cells.Item(699, 11).Value = 600
$ws.Cells.Item(699, 12).Value = 700
$ws.Cells.Item(699, 13).Value = 650
$ws.Cells.Item(699, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(699, 15).Value = "Perú"
$ws.Cells.Item(699, 16).Value = 650
$ws.Cells.Item(699, 17).Value = 1
$ws.Cells.Item(699, 18).Value = "Hortaliza"

# Populate new row 700
$ws.Cells.Item(700, 1).Value = 9
$ws.Cells.Item(700, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(700, 3).Value = "Metropolitana"
$ws.Cells.Item(700, 4).Value = 44509
$ws.Cells.Item(700, 5).Value = 13
$ws.Cells.Item(700, 6).Value = 100112045
$ws.Cells.Item(700, 7).Value = "Zapallo"
$ws.Cells.Item(700, 8).Value = "Camote"
$ws.Cells.Item(700, 9).Value = "2a nueva(o)"
$ws.Cells.Item(700, 10).Value = 61
$ws.Cells.Item(700, 11).Value = 450
$ws.Cells.Item(700, 12).Value = 500
$ws.Cells.Item(700, 13).Value = 475
$ws.Cells.Item(700, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(700, 15).Value = "Perú"
$ws.Cells.Item(700, 16).Value = 475
$ws.Cells.Item(700, 17).Value = 1
$ws.Cells.Item(700, 18).Value = "Hortaliza"
